# Finish notebook analysis: insert a new data row (Content sample collected
# 2021-11-10-17.00) right after the existing row 13, which shifts every
# subsequent data row down by one, and append one more new data row
# (Mobile_App sample collected 2021-11-10-17.28) at the end of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before row 14; this pushes old rows 14-22 down to 15-23.
$ws.Rows.Item(14).Insert()

# Populate the newly inserted row 14 with the "Content" sample.
$ws.Cells.Item(14, 1).Value = 0.7032967032967034
$ws.Cells.Item(14, 2).Value = 0.02197802197802198
$ws.Cells.Item(14, 3).Value = 0.978021978021978
$ws.Cells.Item(14, 4).Value = 0.978021978021978
$ws.Cells.Item(14, 5).Value = 0.978021978021978
$ws.Cells.Item(14, 6).Value = 0.978021978021978
$ws.Cells.Item(14, 7).Value = 0.281043956043956
$ws.Cells.Item(14, 8).Value = 0.4890109890109889
$ws.Cells.Item(14, 9).Value = 0.770054945054945
$ws.Cells.Item(14, 10).Value = 2038
$ws.Cells.Item(14, 11).Value = "Content"
$ws.Cells.Item(14, 12).Value = "2021-11-10-17.00"
$ws.Cells.Item(14, 13).Value = "2021-11-10-17.00"

# Append a new row 24 with the final "Mobile_App" sample.
$ws.Cells.Item(24, 1).Value = 0.5444444444444444
$ws.Cells.Item(24, 2).Value = 0
$ws.Cells.Item(24, 3).Value = 1
$ws.Cells.Item(24, 4).Value = 1
$ws.Cells.Item(24, 5).Value = 1
$ws.Cells.Item(24, 6).Value = 1
$ws.Cells.Item(24, 7).Value = 0.2548333333333334
$ws.Cells.Item(24, 8).Value = 0.5
$ws.Cells.Item(24, 9).Value = 0.7548333333333334
$ws.Cells.Item(24, 10).Value = 2439
$ws.Cells.Item(24, 11).Value = "Mobile_App"
$ws.Cells.Item(24, 12).Value = "2021-11-10-17.28"
$ws.Cells.Item(24, 13).Value = "2021-11-10-17.28"
